# Apply "updated structure of analysis notebooks and result exports" edit.
#
# Semantics of the change (reverse-engineered from the OOXML diff):
#  - The shared-string table order for "pre"/"post" was swapped, and
#    correspondingly the raw data rows that were tied to "pre"/"post" labels
#    in the `normality` and `equal_var` sheets were swapped back (row pairs)
#    so the figures stay attached to the correct label after the
#    source notebook re-ran / re-exported with a different condition order.
#  - Across every results sheet, the previously full-precision floating
#    point statistics (W, pval, SS, F, p-unc, np2, T, p-unc, hedges, etc.)
#    were re-exported rounded to 4 decimal places.
#  - A couple of boolean "normal" flags flipped along with their row's data
#    swap in the `normality` sheet.
#
# Below we simply (re)write each affected cell with its final value.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "normality": pre/post row pairs swapped + values rounded to 4dp
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("normality")
$ws.Range("B3").Value = "post"
$ws.Range("D3").Value = 0.9145
$ws.Range("E3").Value = 0.2436
$ws.Range("D4").Value = 0.955
$ws.Range("E4").Value = 0.676
$ws.Range("B5").Value = "pre"
$ws.Range("D5").Value = 0.9478
$ws.Range("E5").Value = 0.6051
$ws.Range("D6").Value = 0.8975
$ws.Range("E6").Value = 0.1234
$ws.Range("B7").Value = "post"
$ws.Range("D7").Value = 0.9382
$ws.Range("E7").Value = 0.4755
$ws.Range("F7").Value = $true
$ws.Range("D8").Value = 0.9745
$ws.Range("E8").Value = 0.9412
$ws.Range("B9").Value = "pre"
$ws.Range("D9").Value = 0.824
$ws.Range("E9").Value = 0.0178
$ws.Range("F9").Value = $false
$ws.Range("D10").Value = 0.9656
$ws.Range("E10").Value = 0.8362
$ws.Range("B11").Value = "post"
$ws.Range("D11").Value = 0.9593
$ws.Range("E11").Value = 0.7734
$ws.Range("D12").Value = 0.9194
$ws.Range("E12").Value = 0.2464
$ws.Range("B13").Value = "pre"
$ws.Range("D13").Value = 0.9498
$ws.Range("E13").Value = 0.6347
$ws.Range("D14").Value = 0.9096
$ws.Range("E14").Value = 0.1809

# ---------------------------------------------------------------------
# Sheet "equal_var": pre/post row pairs swapped + values rounded to 4dp
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("equal_var")
$ws.Range("B3").Value = "post"
$ws.Range("D3").Value = 0.3335
$ws.Range("E3").Value = 0.5692
$ws.Range("B4").Value = "pre"
$ws.Range("D4").Value = 1.6634
$ws.Range("E4").Value = 0.21
$ws.Range("B5").Value = "post"
$ws.Range("D5").Value = 0.0203
$ws.Range("E5").Value = 0.8879
$ws.Range("B6").Value = "pre"
$ws.Range("D6").Value = 0.7395
$ws.Range("E6").Value = 0.3987
$ws.Range("B7").Value = "post"
$ws.Range("D7").Value = 0.3599
$ws.Range("E7").Value = 0.5544
$ws.Range("B8").Value = "pre"
$ws.Range("D8").Value = 0.1259
$ws.Range("E8").Value = 0.726

# ---------------------------------------------------------------------
# Sheet "mixed_anova": values rounded to 4dp (no row reordering here)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("mixed_anova")
$ws.Range("H3").Value = 2.3558
$ws.Range("I3").Value = 0.1385
$ws.Range("J3").Value = 0.0929
$ws.Range("D4").Value = 131.22
$ws.Range("G4").Value = 131.22
$ws.Range("H4").Value = 9.5185
$ws.Range("I4").Value = 0.0052
$ws.Range("J4").Value = 0.2927
$ws.Range("D5").Value = 1.2063
$ws.Range("G5").Value = 1.2063
$ws.Range("H5").Value = 0.0875
$ws.Range("I5").Value = 0.77
$ws.Range("J5").Value = 0.0038
$ws.Range("D6").Value = 9.8371
$ws.Range("G6").Value = 9.8371
$ws.Range("H6").Value = 0.1762
$ws.Range("I6").Value = 0.6785
$ws.Range("J6").Value = 0.0076
$ws.Range("H7").Value = 56.4902
$ws.Range("I7").Value = 0.0
$ws.Range("J7").Value = 0.7107
$ws.Range("D8").Value = 16.157
$ws.Range("G8").Value = 16.157
$ws.Range("H8").Value = 0.6498
$ws.Range("I8").Value = 0.4284
$ws.Range("J8").Value = 0.0275
$ws.Range("D9").Value = 2.3713
$ws.Range("G9").Value = 2.3713
$ws.Range("H9").Value = 0.0538
$ws.Range("I9").Value = 0.8186
$ws.Range("J9").Value = 0.0023
$ws.Range("D10").Value = 2060.82
$ws.Range("G10").Value = 2060.82
$ws.Range("H10").Value = 65.7966
$ws.Range("I10").Value = 0.0
$ws.Range("J10").Value = 0.741
$ws.Range("D11").Value = 0.2954
$ws.Range("G11").Value = 0.2954
$ws.Range("H11").Value = 0.0094
$ws.Range("I11").Value = 0.9235
$ws.Range("J11").Value = 0.0004

# ---------------------------------------------------------------------
# Sheet "pairwise_ttests": values rounded to 4dp (labels unchanged)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("pairwise_ttests")
$ws.Range("I3").Value = -3.1456
$ws.Range("L3").Value = 0.0044
$ws.Range("N3").Value = -0.5846
$ws.Range("I4").Value = 1.5485
$ws.Range("J4").Value = 22.6175
$ws.Range("L4").Value = 0.1354
$ws.Range("N4").Value = 0.5942
$ws.Range("I5").Value = 1.4334
$ws.Range("J5").Value = 21.6593
$ws.Range("L5").Value = 0.166
$ws.Range("N5").Value = 0.5588
$ws.Range("I6").Value = 1.2469
$ws.Range("J6").Value = 21.3984
$ws.Range("L6").Value = 0.2259
$ws.Range("N6").Value = 0.4757
$ws.Range("I7").Value = -7.5714
$ws.Range("L7").Value = 0.0
$ws.Range("N7").Value = -1.6664
$ws.Range("I8").Value = -0.4152
$ws.Range("J8").Value = 20.533
$ws.Range("L8").Value = 0.6823
$ws.Range("N8").Value = -0.1625
$ws.Range("I9").Value = 0.0868
$ws.Range("J9").Value = 22.969
$ws.Range("L9").Value = 0.9316
$ws.Range("N9").Value = 0.0335
$ws.Range("I10").Value = -0.9294
$ws.Range("J10").Value = 18.7651
$ws.Range("L10").Value = 0.3645
$ws.Range("N10").Value = -0.3658
$ws.Range("I11").Value = -8.2843
$ws.Range("L11").Value = 0.0
$ws.Range("N11").Value = -2.1015
$ws.Range("I12").Value = -0.2339
$ws.Range("J12").Value = 22.7069
$ws.Range("L12").Value = 0.8171
$ws.Range("N12").Value = -0.0898
$ws.Range("I13").Value = -0.1891
$ws.Range("J13").Value = 22.9381
$ws.Range("L13").Value = 0.8517
$ws.Range("N13").Value = -0.0728
$ws.Range("I14").Value = -0.1884
$ws.Range("J14").Value = 22.8679
$ws.Range("L14").Value = 0.8522
$ws.Range("N14").Value = -0.0729

Write-Output "done"
